# Adding new case (row 46), text only - no msg pic.
# Mirrors commit: "adding new case, text only no msg pic"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unrelated cleanup that shipped in the same commit: row 7's motivation
# cell (F7, "financial") was cleared out. ---
$ws.Range("F7").ClearContents()

# --- New row 46: id 45, same template as row 45 but a different case. ---
$ws.Range("A46").Value = 45
$ws.Range("B46").Value = "msg"

# Date: copy value + the date number-format from C45 so we reuse its style
# (08/08/2021, same "added" date as the previous entry) instead of minting a
# brand-new style entry.
$ws.Range("C46").Value = 44416
$ws.Range("C45").Copy()
$ws.Range("C46").PasteSpecial(-4122)

$ws.Range("D46").Value = "MCAST"
$ws.Range("E46").Value = "shortened"
$ws.Range("F46").Value = "financial"
$ws.Range("G46").Value = "mt"
$ws.Range("H46").Value = "no"

# Entity (J) filled before description (I) so new shared strings are
# interned in the same order as the source commit.
$ws.Range("J46").Value = "GO,Melita"
$ws.Range("I46").Value = "refund from telecoms"
$ws.Range("K46").Value = "redirects to https://doctorbrew.pl/wp-admin/user/-/"

# Move the view/selection down to the newly added row.
$excel.ActiveWindow.ScrollRow = 32
$ws.Range("E46").Select()
